$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-6 currently hold (Activity / ScheduledResource / Plant / ProductionDivision):
#   Row4: Pos Cut / HC Digital Knife / HC Digital Knife / HC Digital Knife
#   Row5: Proof Approval / HC Proof Approval / HC Proof  Approval / HC Proof Approval
#   Row6: Digital Print F/B 4x4 / HC NexPress 1-4c / HC NexPress 1-4c / HC NexPress 1-4c
#
# Target: row 6's content moves up to row 4, pushing the other two rows down one.
#   Row4: Digital Print F/B 4x4 / HC NexPress 1-4c / HC NexPress 1-4c / HC NexPress 1-4c
#   Row5: Pos Cut / HC Digital Knife / HC Digital Knife / HC Digital Knife
#   Row6: Proof Approval / HC Proof Approval / HC Proof  Approval / HC Proof Approval

$ws.Range("B4").Value = "Digital Print F/B 4x4"
$ws.Range("D4").Value = "HC NexPress 1-4c"
$ws.Range("L4").Value = "HC NexPress 1-4c"
$ws.Range("M4").Value = "HC NexPress 1-4c"

$ws.Range("B5").Value = "Pos Cut"
$ws.Range("D5").Value = "HC Digital Knife"
$ws.Range("L5").Value = "HC Digital Knife"
$ws.Range("M5").Value = "HC Digital Knife"

$ws.Range("B6").Value = "Proof Approval"
$ws.Range("D6").Value = "HC Proof Approval"
$ws.Range("L6").Value = "HC Proof  Approval"
$ws.Range("M6").Value = "HC Proof Approval"
